# Apply the edits described by the commit:
#  1. Date placeholder text on the slide master + all 11 custom layouts:
#     8/24/2010 -> 8/25/2010
#  2. Slide 1 ("Flowchart: Multidocument 19"): collections -> Collections
#  3. Slide 2 ("Cloud Callout 6"): Tile -> Title (2nd paragraph)
#  4. Slide 3 ("TextBox 5"): merge runs "(2) " + "Module " into "(2) Module "

$p = $ppt.ActivePresentation

# --- 1. Date placeholders -------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "8/24/2010") {
            $shp.TextFrame.TextRange.Text = "8/25/2010"
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "8/24/2010") {
                $shp.TextFrame.TextRange.Text = "8/25/2010"
            }
        }
    }
}

# Helper: find a shape on a slide by its name.
function Find-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $slide.Shapes.Item(1)
}

# --- 2. Slide 1: collections -> Collections -------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = Find-ShapeByName $s1 "Flowchart: Multidocument 19"
$shp1.TextFrame.TextRange.Text = "Collections"

# --- 3. Slide 2: Tile -> Title --------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = Find-ShapeByName $s2 "Cloud Callout 6"
$tr2 = $shp2.TextFrame.TextRange
$found = $tr2.Find("Tile", 0)
$found.Text = "Title"

# --- 4. Slide 3: merge "(2) " + "Module " runs ----------------------------
$s3 = $p.Slides.Item(3)
$shp3 = Find-ShapeByName $s3 "TextBox 5"
$tr3 = $shp3.TextFrame.TextRange
# "Module " is characters 5-11 of the run text; delete it, then append it
# back onto the end of the "(2) " run so the two runs merge into one.
$moduleRange = $tr3.Characters(5, 7)
[void]$moduleRange.Delete()
$tr3b = $shp3.TextFrame.TextRange
$firstRun = $tr3b.Characters(1, 4)
[void]$firstRun.InsertAfter("Module ")
